$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new E value, new F value) updates, derived from the countdown logic:
# each rows "remaining days" (E) decrements by 1 per day; when it would hit 0
# (i.e. was 1), it resets to 10 and the start date (F) advances by 10 days.
# Row 36 is excluded (its F value is a malformed date and was left untouched).
$updates = @(
    @{Row=2; E=10; F=20251215}
    @{Row=3; E=10; F=20251215}
    @{Row=4; E=10; F=20251215}
    @{Row=5; E=6; F=20251215}
    @{Row=6; E=10; F=20251215}
    @{Row=7; E=6; F=20251215}
    @{Row=8; E=10; F=20251215}
    @{Row=9; E=6; F=20251215}
    @{Row=10; E=3; F=20251215}
    @{Row=11; E=10; F=20251215}
    @{Row=12; E=6; F=20251215}
    @{Row=13; E=10; F=20251215}
    @{Row=14; E=10; F=20251215}
    @{Row=15; E=10; F=20251215}
    @{Row=16; E=10; F=20251219}
    @{Row=17; E=6; F=20251215}
    @{Row=18; E=9; F=20251218}
    @{Row=19; E=9; F=20251218}
    @{Row=20; E=9; F=20251218}
    @{Row=21; E=9; F=20251218}
    @{Row=22; E=6; F=20251215}
    @{Row=23; E=6; F=20251215}
    @{Row=24; E=6; F=20251215}
    @{Row=25; E=6; F=20251215}
    @{Row=26; E=6; F=20251215}
    @{Row=27; E=4; F=20251216}
    @{Row=28; E=9; F=20251218}
    @{Row=29; E=9; F=20251218}
    @{Row=30; E=9; F=20251218}
    @{Row=31; E=9; F=20251218}
    @{Row=32; E=9; F=20251218}
    @{Row=33; E=9; F=20251218}
    @{Row=34; E=9; F=20251218}
    @{Row=35; E=9; F=20251218}
    @{Row=37; E=9; F=20251218}
    @{Row=38; E=9; F=20251218}
    @{Row=39; E=9; F=20251218}
    @{Row=40; E=3; F=20251215}
    @{Row=41; E=3; F=20251215}
    @{Row=42; E=9; F=20251218}
    @{Row=43; E=6; F=20251215}
    @{Row=44; E=3; F=20251215}
    @{Row=45; E=6; F=20251215}
    @{Row=46; E=3; F=20251215}
    @{Row=47; E=9; F=20251218}
    @{Row=48; E=3; F=20251215}
    @{Row=49; E=4; F=20251216}
    @{Row=50; E=4; F=20251213}
    @{Row=51; E=4; F=20251213}
    @{Row=52; E=4; F=20251213}
    @{Row=53; E=4; F=20251213}
    @{Row=54; E=4; F=20251213}
    @{Row=55; E=4; F=20251213}
    @{Row=56; E=4; F=20251213}
    @{Row=57; E=4; F=20251213}
    @{Row=58; E=8; F=20251217}
    @{Row=59; E=8; F=20251217}
    @{Row=60; E=8; F=20251217}
    @{Row=61; E=4; F=20251216}
    @{Row=62; E=8; F=20251217}
    @{Row=63; E=8; F=20251217}
    @{Row=64; E=8; F=20251217}
    @{Row=65; E=9; F=20251218}
    @{Row=66; E=9; F=20251218}
    @{Row=67; E=9; F=20251218}
    @{Row=68; E=9; F=20251218}
    @{Row=69; E=9; F=20251218}
    @{Row=70; E=10; F=20251219}
    @{Row=71; E=10; F=20251219}
    @{Row=72; E=10; F=20251219}
    @{Row=73; E=10; F=20251219}
    @{Row=74; E=10; F=20251219}
    @{Row=75; E=10; F=20251219}
    @{Row=76; E=10; F=20251219}
    @{Row=77; E=3; F=20251212}
    @{Row=78; E=3; F=20251212}
    @{Row=79; E=3; F=20251212}
    @{Row=80; E=3; F=20251212}
    @{Row=81; E=3; F=20251212}
    @{Row=82; E=3; F=20251212}
    @{Row=83; E=3; F=20251212}
    @{Row=84; E=3; F=20251212}
    @{Row=85; E=3; F=20251212}
    @{Row=86; E=3; F=20251212}
    @{Row=87; E=3; F=20251215}
    @{Row=88; E=3; F=20251215}
    @{Row=89; E=3; F=20251215}
    @{Row=90; E=3; F=20251215}
    @{Row=91; E=6; F=20251215}
    @{Row=92; E=3; F=20251215}
    @{Row=93; E=3; F=20251212}
    @{Row=94; E=6; F=20251218}
    @{Row=95; E=2; F=20251211}
    @{Row=96; E=10; F=20251219}
    @{Row=97; E=10; F=20251219}
    @{Row=98; E=10; F=20251219}
    @{Row=99; E=10; F=20251219}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}

